# Add AD.SEC.001.FON.99 Ouvrir session OU reconnection pour moi :-)
#
# - A2 now holds the new shared string "AD.SEC.001.FON.99" (replacing
#   "AD.SEC.001.FON.01").
# - Row 3's test-case columns A and B (the "RO.ACT.003HAB.SRM" / 20 pair)
#   are removed.
# - The active selection moves from B4 to B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AD.SEC.001.FON.99"

$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()

$ws.Range("B12").Select()
